$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '37.200.50'
$ws.Range("E2").Value = '  -0.33%  '
# Row 3
$ws.Range("D3").Value = '2.027.75'
$ws.Range("E3").Value = '  -0.97%  '
# Row 4
$ws.Range("E4").Value = '  -0.28%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.22'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.10%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.604'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.59%  '
# Row 7
$ws.Range("E7").Value = '  -0.06%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '55.31'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.54%  '
# Row 9
$ws.Range("E9").Value = '  -1.38%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0787'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.27%  '
# Row 11
$ws.Range("E11").Value = '  -5.17%  '
# Row 12
$ws.Range("D12").Value = '2.321.72'
$ws.Range("E12").Value = '  -1.33%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.30'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.71%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.33'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.10%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.745'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.37%  '
# Row 16
$ws.Range("E16").Value = '  -1.86%  '
# Row 17
$ws.Range("D17").Value = '2.025.67'
$ws.Range("E17").Value = '  -1.03%  '
# Row 18
$ws.Range("D18").Value = '37.142.98'
$ws.Range("E18").Value = '  -0.19%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.54'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +7.61%  '
# Row 20
$ws.Range("E20").Value = '  -0.61%  '
# Row 21
$ws.Range("E21").Value = '  -1.14%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '223.84'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.03%  '
# Row 23
$ws.Range("E23").Value = '  +0.15%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.43'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.60%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.19'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.33%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.26'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.01%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.24'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.79%  '
# Row 28
$ws.Range("E28").Value = '  -1.10%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.77'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.36%  '
# Row 30
$ws.Range("E30").Value = '  -1.56%  '
# Row 31
$ws.Range("E31").Value = '  -1.20%  '
# Row 32
$ws.Range("E32").Value = '  -0.08%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0616'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.11%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.50'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.94%  '
# Row 35
$ws.Range("E35").Value = '  -5.00%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.86'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.59%  '
# Row 37
$ws.Range("E37").Value = '  -0.23%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.56'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.02%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.13'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.02%  '
# Row 40
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0216'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.13%  '
# Row 41
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '1.471.30'
$ws.Range("E41").Value = '  -0.80%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '96.08'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.58%  '
# Row 43
$ws.Range("B43").Value = 'HuobiToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.82'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.64%  '
# Row 44
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.46'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.92%  '
# Row 45
$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0912'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.41%  '
# Row 46
$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.14'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.19%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.27'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.93%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.02'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.58%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.94'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.68%  '
# Row 50
$ws.Range("E50").Value = '  -7.59%  '
# Row 51
$ws.Range("D51").Value = '2.206.29'
$ws.Range("E51").Value = '  -1.43%  '
